$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.561.33"
$ws.Range("E2").Value = "  +3.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.587.50"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.95%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.00"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.32"
$ws.Range("E8").Value = "  +6.39%  "
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0601"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0885"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.814.29"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.594.73"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.541.38"
$ws.Range("E16").Value = "  +3.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.08"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.82"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0706"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.05"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.33"
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.67"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.25"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.55"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0470"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.398.21"
$ws.Range("E34").Value = "  -3.94%  "
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  -9.99%  "
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("E38").Value = "  +10.51%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.812"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  +0.95%  "
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.982"
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.01"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.724.28"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.30"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0104"
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "41.54"
$ws.Range("E51").Value = "  +12.75%  "
